$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.137.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.047.76"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'248.24"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").Value = "'0.664"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D8").Value = "'56.99"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'0.0783"
$ws.Range("E10").Value = "  -2.33%  "
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "'16.23"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'0.895"
$ws.Range("E13").Value = "  +10.30%  "
$ws.Range("D14").Value = "2.345.01"
$ws.Range("D15").Value = "'5.72"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "2.048.42"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "'18.51"
$ws.Range("E17").Value = "  +12.64%  "
$ws.Range("D18").Value = "37.145.74"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "'74.47"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("E20").Value = "  -3.45%  "
$ws.Range("D21").Value = "'5.43"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'236.70"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "'9.60"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").Value = "'170.16"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").Value = "'20.14"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").Value = "'4.95"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "'0.0884"
$ws.Range("E34").Value = "  -3.31%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'2.26"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").Value = "'5.28"
$ws.Range("E39").Value = "  +12.72%  "
$ws.Range("E40").Value = "  +8.32%  "
$ws.Range("D41").Value = "'0.0992"
$ws.Range("E41").Value = "  -15.39%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'17.59"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0223"
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "'95.84"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("E46").Value = "  -2.30%  "
$ws.Range("D47").Value = "1.271.28"
$ws.Range("E47").Value = "  -2.31%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "'6.79"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").Value = "2.231.22"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").Value = "'44.41"
$ws.Range("E51").Value = "  +0.24%  "
